$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 'GEwK0Qcc'
$ws.Range("B13").Value = "'11/10/2024"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = '19:00'
$ws.Range("D13").Value = 'VENEZUELA - LIGA FUTVE'
$ws.Range("E13").Value = 'La Guaira'
$ws.Range("F13").Value = 'Rayo Zuliano'
$ws.Range("G13").Value = 1.83
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 4.25
$ws.Range("J13").Value = 2.47
$ws.Range("K13").Value = 1.98
$ws.Range("L13").Value = 4.8
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 6.95
$ws.Range("O13").Value = 1.39
$ws.Range("P13").Value = 2.57
$ws.Range("Q13").Value = 2.12
$ws.Range("R13").Value = 1.57
$ws.Range("S13").Value = 1.47
$ws.Range("T13").Value = 2.32
$ws.Range("U13").Value = 1.98
$ws.Range("V13").Value = 1.65
$ws.Range("W13").Value = 5.9
$ws.Range("X13").Value = 7.8
$ws.Range("Y13").Value = 8.5
$ws.Range("Z13").Value = 15
$ws.Range("AA13").Value = 16
$ws.Range("AB13").Value = 35
$ws.Range("AC13").Value = 7.7
$ws.Range("AD13").Value = 6.4
$ws.Range("AE13").Value = 18
$ws.Range("AF13").Value = 110
$ws.Range("AG13").Value = 9.75
$ws.Range("AH13").Value = 22
$ws.Range("AI13").Value = 15
$ws.Range("AJ13").Value = 70
$ws.Range("AK13").Value = 50
$ws.Range("AL13").Value = 60
$ws.Range("AM13").Value = 201
$ws.Range("AN13").Value = 3.5
$ws.Range("AO13").Value = 9.5
$ws.Range("AP13").Value = 21
$ws.Range("AQ13").Value = 35
$ws.Range("AR13").Value = 80
$ws.Range("AS13").Value = 350
$ws.Range("AT13").Value = 2.27
$ws.Range("AU13").Value = 8
$ws.Range("AV13").Value = 90
$ws.Range("AW13").Value = 5.8
$ws.Range("AX13").Value = 26
$ws.Range("AY13").Value = 37
$ws.Range("AZ13").Value = 175
$ws.Range("BA13").Value = 250
$ws.Range("BB13").Value = 51
$ws.Range("BC13").Value = 51
$ws.Range("BD13").Value = 51

# Row 14
$ws.Range("A14").Value = '4QWYxN5j'
$ws.Range("B14").Value = "'11/10/2024"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '21:30'
$ws.Range("D14").Value = 'VENEZUELA - LIGA FUTVE'
$ws.Range("E14").Value = 'Zamora'
$ws.Range("F14").Value = 'Monagas'
$ws.Range("G14").Value = 3.3
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 3.8
$ws.Range("K14").Value = 2.12
$ws.Range("L14").Value = 2.6
$ws.Range("M14").Value = 8.300000000000001
$ws.Range("N14").Value = 1.05
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 1.91
$ws.Range("R14").Value = 1.8
$ws.Range("S14").Value = 1.39
$ws.Range("T14").Value = 2.57
$ws.Range("U14").Value = 1.75
$ws.Range("V14").Value = 1.87
$ws.Range("W14").Value = 10
$ws.Range("X14").Value = 17
$ws.Range("Y14").Value = 11.75
$ws.Range("Z14").Value = 45
$ws.Range("AA14").Value = 29
$ws.Range("AB14").Value = 37
$ws.Range("AC14").Value = 10.25
$ws.Range("AD14").Value = 6.8
$ws.Range("AE14").Value = 15
$ws.Range("AF14").Value = 70
$ws.Range("AG14").Value = 7.3
$ws.Range("AH14").Value = 9.5
$ws.Range("AI14").Value = 8.75
$ws.Range("AJ14").Value = 17.5
$ws.Range("AK14").Value = 16
$ws.Range("AL14").Value = 28
$ws.Range("AM14").Value = 600
$ws.Range("AN14").Value = 5.1
$ws.Range("AO14").Value = 18
$ws.Range("AP14").Value = 26
$ws.Range("AQ14").Value = 90
$ws.Range("AR14").Value = 120
$ws.Range("AS14").Value = 350
$ws.Range("AT14").Value = 2.52
$ws.Range("AU14").Value = 7.4
$ws.Range("AV14").Value = 70
$ws.Range("AW14").Value = 3.85
$ws.Range("AX14").Value = 10.25
$ws.Range("AY14").Value = 19.5
$ws.Range("AZ14").Value = 40
$ws.Range("BA14").Value = 75
$ws.Range("BB14").Value = 250
$ws.Range("BC14").Value = 51
$ws.Range("BD14").Value = 51

